# Weekly update: a new price record is inserted as row 103 (pushing the
# existing rows 103:193 down to 104:194, which is why the sheet's last row
# grows from 193 to 194). Populate the freshly inserted row with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103; everything below shifts down by one.
$ws.Rows.Item(103).Insert()

# Fill in the new row 103 with the new record (same market/product as the
# surrounding rows, new date + volume + prices).
$ws.Range("A103").Value = 5
$ws.Range("B103").Value = "Macroferia Regional de Talca"
$ws.Range("C103").Value = "Maule"
$ws.Range("D103").Value = 44512
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = 100112006
$ws.Range("G103").Value = "Repollo"
$ws.Range("H103").Value = "Crespo record"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 4000
$ws.Range("K103").Value = 800
$ws.Range("L103").Value = 800
$ws.Range("M103").Value = 800
$ws.Range("N103").Value = "$/unidad"
$ws.Range("O103").Value = "Región del Maule"
$ws.Range("P103").Value = 800
$ws.Range("Q103").Value = 1
$ws.Range("R103").Value = "Hortaliza"
